# "6. solicitacao_fers.xlsx" - ER atualizado e tabelas mapeadas na pasta 'tabelas novas'
#
# Changes applied:
#  1. B2 ("Tabela Nome novo" column, row for "solicitacao_fers") is updated
#     from the old mapped table name "vacation_solicitations" to a
#     placeholder "XXX" (table name not yet decided / to be filled in).
#  2. B3 ("Novo" column header row, value "Novo") gets its formatting
#     refreshed/re-applied (re-stamping the font), which is what produced
#     the extra cell style in the saved workbook.
#  3. Page setup is configured to Portrait / A4-ish "paperSize 9" printing,
#     matching the newly added <pageSetup> element.
#  4. The active selection is left on B3, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the mapped table name placeholder in B2.
$ws.Range("B2").Value = "XXX"

# 2) Re-apply the font on B3 so the cell gets its own explicit style
#    (keeps the existing blue fill, re-stamps the font face).
$ws.Range("B3").Font.Name = "Calibri"

# 3) Configure the page for printing (Portrait, paper size 9 = A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# 4) Leave the selection on B3.
$ws.Range("B3").Select()
